$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 245 (everything from the old row 245 onward
# shifts down by two rows, old row 294/295 become new 296/297).
$ws.Rows.Item(245).Insert()
$ws.Rows.Item(245).Insert()

# --- New row 245: Angeleno / Primera, Region de O'Higgins ---
$ws.Range("A245").Value = 8
$ws.Range("B245").Value = "Terminal La Palmera de La Serena"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44711
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100103
$ws.Range("H245").Value = "Frutos de hueso (carozo)"
$ws.Range("I245").Value = 100103002
$ws.Range("J245").Value = "Ciruela"
$ws.Range("K245").Value = "Angeleno"
$ws.Range("L245").Value = "Primera"
$ws.Range("M245").Value = 20
$ws.Range("N245").Value = 200000
$ws.Range("O245").Value = 210000
$ws.Range("P245").Value = 205000
$ws.Range("Q245").Value = "$/bins (450 kilos)"
$ws.Range("R245").Value = "Región de O'Higgins"
$ws.Range("S245").Value = 456
$ws.Range("T245").Value = 450

# --- New row 246: Angeleno / Segunda, Region de O'Higgins ---
$ws.Range("A246").Value = 8
$ws.Range("B246").Value = "Terminal La Palmera de La Serena"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = 44711
$ws.Range("E246").Value = 4
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100103
$ws.Range("H246").Value = "Frutos de hueso (carozo)"
$ws.Range("I246").Value = 100103002
$ws.Range("J246").Value = "Ciruela"
$ws.Range("K246").Value = "Angeleno"
$ws.Range("L246").Value = "Segunda"
$ws.Range("M246").Value = 16
$ws.Range("N246").Value = 170000
$ws.Range("O246").Value = 180000
$ws.Range("P246").Value = 175000
$ws.Range("Q246").Value = "$/bins (450 kilos)"
$ws.Range("R246").Value = "Región de O'Higgins"
$ws.Range("S246").Value = 389
$ws.Range("T246").Value = 450
